$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 66.09241856096124

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 14).Value = $newValue
}
